$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Consolidate the title's text runs ("A" + " " + "slide") into a single run.
# Setting the same logical text is a no-op for the engine's run model, so we
# first flip the text to a throwaway value to force a rebuild of the runs,
# then set the final desired text.
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "~~tmp~~"
$title.Text = "A slide"

# Consolidate the caption textbox's many single-word runs into one run.
$caption = $s.Shapes.Item(4).TextFrame.TextRange
$caption.Text = "~~tmp~~"
$caption.Text = "Just an image on this side"
